# Add 2022-Q4 data
# 1) Insert a new worksheet "2022-Q4" right before the existing "2022-Q3" sheet
#    and populate it with the fund-holding detail for that quarter.
# 2) Insert a new row at the top of the "总计" (totals) summary sheet with the
#    2022-Q4 aggregate numbers (count of holdings + total market value).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Part 1: create the "2022-Q4" worksheet
# ---------------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$newSheet = $wb.Worksheets.Add($q3Sheet)
$newSheet.Name = "2022-Q4"

# Header row
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $cell = $newSheet.Cells.Item(1, $col)
    $cell.Value = $headers[$col - 2]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# Data rows: index, 基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名(number)
$rows = @(
    @("0", "217024", "招商安盈债券A", "40.95", "20.32", "1.40", "0.5733", 5),
    @("1", "014887", "招商安福1年定期开放债券", "17.78", "33.59", "1.23", "0.2187", 8),
    @("2", "010430", "招商安阳债券A", "16.90", "20.35", "0.86", "0.1453", 7),
    @("3", "016513", "招商安嘉债券", "16.27", "20.17", "0.84", "0.1367", 6),
    @("4", "009658", "汇丰晋信中小盘低波动策略股票A", "0.85", "92.42", "1.47", "0.0125", 4),
    @("5", "010431", "招商安阳债券C", "0.11", "20.35", "0.86", "0.0009", 7),
    @("6", "009775", "汇丰晋信中小盘低波动策略股票C", "0.04", "92.42", "1.47", "0.0006", 4),
    @("7", "012233", "招商安盈债券C", "0.01", "20.32", "1.40", "0.0001", 5)
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r, 1).Value = [double]$row[0]
    $newSheet.Cells.Item($r, 2).Value = "'" + $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]

    $dCell = $newSheet.Cells.Item($r, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $row[3]

    $eCell = $newSheet.Cells.Item($r, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $row[4]

    $fCell = $newSheet.Cells.Item($r, 6)
    $fCell.NumberFormat = "@"
    $fCell.Value = $row[5]

    $gCell = $newSheet.Cells.Item($r, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value = $row[6]

    $newSheet.Cells.Item($r, 8).Value = $row[7]

    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Part 2: insert the 2022-Q4 row into the "总计" summary sheet
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows("2:2").Insert()
$totalSheet.Range("A2:D2").ClearFormats()

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 8
$totalSheet.Range("D2").Value = 1.09

# ---------------------------------------------------------------------------
# Restore the originally-active tab (last sheet, "2021-Q3") as the selected one
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q3").Activate()
